$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column B. This pushes the existing
# "Jun_13" column (old B) to D and the existing "Jun_10" column (old C) to E,
# carrying their values/styles along (e.g. the highlighted cell in row 18).
$ws.Columns("B:C").Insert()

# New header row: two new date columns "Jun_17" / "Jun_15" in B1 / C1.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the two new data columns (B, C) on every data row with the same
# placeholder analyst-rank value ("UN") that used to live in column B
# before the insert. (The insert cleared B; D already has "UN" carried
# over automatically from the old column B, and E has the old column C's
# value, so only B and C need to be (re)populated here.)
$lastRow = 27
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Match the new columns' width to the existing formatted column (8 chars).
$ws.Columns("C:C").ColumnWidth = 7.166666666666667
$ws.Columns("D:D").ColumnWidth = 7.166666666666667
$ws.Columns("E:E").ColumnWidth = 7.166666666666667
